# Commit: "Added restart game functionality."
#
# The scoreboard on Scores gets a fresh round appended below the existing
# data (rows 2-40): two more "1" round-marker rows reusing the existing
# shared string, and two brand-new players, VITALIJUS and LAJA, each
# scored at 200 like the rest of the recent rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "'1"
$ws.Range("A41").Style = $ws.Range("A1").Style
$ws.Range("B41").Value = 200

$ws.Range("A42").Value = "VITALIJUS"
$ws.Range("B42").Value = 200

$ws.Range("A43").Value = "'1"
$ws.Range("A43").Style = $ws.Range("A1").Style
$ws.Range("B43").Value = 200

$ws.Range("A44").Value = "LAJA"
$ws.Range("B44").Value = 200
